$wb = $excel.ActiveWorkbook

# Sheet 1: 식당판매 (Restaurant sales)
$ws1 = $wb.Worksheets.Item("식당판매")
$ws1.Range("C6").Value = 2
$ws1.Range("C7").ClearContents()
$ws1.Range("C8").ClearContents()

# Sheet 2: 매점판매 (Shop sales)
$ws2 = $wb.Worksheets.Item("매점판매")
$ws2.Range("C2").ClearContents()
$ws2.Range("C3").ClearContents()
$ws2.Range("C5").Value = 2

# Sheet 3: 장의용품 (Funeral supplies)
$ws3 = $wb.Worksheets.Item("장의용품")
$ws3.Range("C2").ClearContents()
$ws3.Range("C7").ClearContents()

# Sheet 4: 상복 (Mourning clothes)
$ws4 = $wb.Worksheets.Item("상복")
$ws4.Range("C3").ClearContents()
$ws4.Range("C6").ClearContents()
$ws4.Range("C7").Value = 5
$ws4.Range("C8").ClearContents()
$ws4.Range("C11").Value = 3

# Sheet 5: 기타 (Other)
$ws5 = $wb.Worksheets.Item("기타")
$ws5.Range("C2").Value = 3
$ws5.Range("C3").Value = 1
$ws5.Range("C5").Value = 1
$ws5.Range("C7").Value = 9
$ws5.Range("C8").Value = 44
$ws5.Range("C9").Value = 5
$ws5.Range("C10").Value = 121
$ws5.Range("C11").Value = 0
$ws5.Range("C12").Value = 2
$ws5.Range("C13").Value = 70
$ws5.Range("C14").Value = 21
$ws5.Range("C15").Value = 23
